$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates. Values that are ambiguous with numbers (single
# decimal point) are written with a leading apostrophe and then restyled
# back to "Normal" so Excel keeps them as plain text (matching the source
# workbook's inlineStr cells) instead of silently coercing them to floats.
$dUpdates = @(
    @{Row=2;  Value='26.889.61'; Quote=$false},
    @{Row=3;  Value='1.564.92';  Quote=$false},
    @{Row=5;  Value='205.92';    Quote=$true},
    @{Row=6;  Value='0.486';     Quote=$true},
    @{Row=8;  Value='21.82';     Quote=$true},
    @{Row=11; Value='0.0865';    Quote=$true},
    @{Row=12; Value='1.786.06';  Quote=$false},
    @{Row=13; Value='1.573.79';  Quote=$false},
    @{Row=16; Value='26.893.07'; Quote=$false},
    @{Row=17; Value='61.37';     Quote=$true},
    @{Row=18; Value='215.31';    Quote=$true},
    @{Row=19; Value='7.38';      Quote=$true},
    @{Row=22; Value='4.14';      Quote=$true},
    @{Row=25; Value='154.32';    Quote=$true},
    @{Row=26; Value='6.69';      Quote=$true},
    @{Row=27; Value='14.96';     Quote=$true},
    @{Row=32; Value='3.16';      Quote=$true},
    @{Row=33; Value='1.392.55';  Quote=$false},
    @{Row=37; Value='0.923';     Quote=$true},
    @{Row=40; Value='0.813';     Quote=$true},
    @{Row=42; Value='0.991';     Quote=$true},
    @{Row=43; Value='5.51';      Quote=$true},
    @{Row=44; Value='1.80';      Quote=$true},
    @{Row=46; Value='63.85';     Quote=$true},
    @{Row=47; Value='1.700.93';  Quote=$false},
    @{Row=48; Value='86.79';     Quote=$true},
    @{Row=49; Value='0.0₇0983';  Quote=$false},
    @{Row=50; Value='0.0503';    Quote=$true},
    @{Row=51; Value='0.0952';    Quote=$true}
)

foreach ($u in $dUpdates) {
    $cell = $ws.Cells.Item($u.Row, 4)
    if ($u.Quote) {
        $cell.Value = "'" + $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}

# Column E (Volume(1h)) updates - always plain text (contains '%' and
# padding spaces so Excel never mistakes them for numbers).
$eUpdates = @(
    @{Row=2;  Value='  -0.94%  '},
    @{Row=3;  Value='  +0.36%  '},
    @{Row=4;  Value='  -0.11%  '},
    @{Row=5;  Value='  -0.38%  '},
    @{Row=6;  Value='  -0.87%  '},
    @{Row=7;  Value='  -0.07%  '},
    @{Row=8;  Value='  -1.41%  '},
    @{Row=9;  Value='  -0.44%  '},
    @{Row=10; Value='  -0.95%  '},
    @{Row=11; Value='  +0.34%  '},
    @{Row=12; Value='  +0.24%  '},
    @{Row=13; Value='  +0.89%  '},
    @{Row=14; Value='  -0.87%  '},
    @{Row=16; Value='  -0.94%  '},
    @{Row=17; Value='  -2.51%  '},
    @{Row=18; Value='  -0.08%  '},
    @{Row=19; Value='  +2.19%  '},
    @{Row=20; Value='  -0.62%  '},
    @{Row=21; Value='  -0.14%  '},
    @{Row=22; Value='  +0.55%  '},
    @{Row=23; Value='  -1.53%  '},
    @{Row=24; Value='  +0.88%  '},
    @{Row=25; Value='  +1.79%  '},
    @{Row=26; Value='  +1.81%  '},
    @{Row=27; Value='  +0.31%  '},
    @{Row=28; Value='  -0.15%  '},
    @{Row=29; Value='  -0.70%  '},
    @{Row=30; Value='  +0.99%  '},
    @{Row=31; Value='  -3.38%  '},
    @{Row=32; Value='  -0.04%  '},
    @{Row=33; Value='  +0.87%  '},
    @{Row=34; Value='  -0.24%  '},
    @{Row=35; Value='  -1.14%  '},
    @{Row=36; Value='  -0.43%  '},
    @{Row=37; Value='  -2.31%  '},
    @{Row=38; Value='  -0.46%  '},
    @{Row=39; Value='  +3.73%  '},
    @{Row=40; Value='  +0.29%  '},
    @{Row=41; Value='  -0.10%  '},
    @{Row=42; Value='  +0.22%  '},
    @{Row=43; Value='  +5.56%  '},
    @{Row=44; Value='  +0.48%  '},
    @{Row=46; Value='  +0.81%  '},
    @{Row=47; Value='  +0.42%  '},
    @{Row=48; Value='  +1.51%  '},
    @{Row=49; Value='  +0.04%  '},
    @{Row=50; Value='  +2.42%  '},
    @{Row=51; Value='  +1.02%  '}
)

foreach ($u in $eUpdates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.Value
}
